$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 258, shifting existing rows 258-264 down to 259-265
$ws.Rows.Item(258).Insert()

# Populate the new row 258 with the latest weekly price entry
$ws.Cells.Item(258, 1).Value = 11
$ws.Cells.Item(258, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(258, 3).Value = "Bíobío"
$ws.Cells.Item(258, 4).Value = 45267
$ws.Cells.Item(258, 5).Value = 8
$ws.Cells.Item(258, 6).Value = "Fruta"
$ws.Cells.Item(258, 7).Value = 100109
$ws.Cells.Item(258, 8).Value = "Uva"
$ws.Cells.Item(258, 9).Value = 100109001
$ws.Cells.Item(258, 10).Value = "Uva"
$ws.Cells.Item(258, 11).Value = "Red Globe"
$ws.Cells.Item(258, 12).Value = "Primera"
$ws.Cells.Item(258, 13).Value = 150
$ws.Cells.Item(258, 14).Value = 16000
$ws.Cells.Item(258, 15).Value = 16000
$ws.Cells.Item(258, 16).Value = 16000
$ws.Cells.Item(258, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(258, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(258, 19).Value = 2000
$ws.Cells.Item(258, 20).Value = 8
